# Auto-generated script applying the Kujata_Profits market-data refresh
# Updates numeric price/profit columns (H-N) on affected rows across sheets
# ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, matching the scheduled-runner commit.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1403.0769
$ws.Range("I18").Value = 1224
$ws.Range("K18").Value = 1224
$ws.Range("M18").Value = -940
$ws.Range("H19").Value = 350.23077
$ws.Range("I19").Value = 275.16666
$ws.Range("J19").Value = 414.57144
$ws.Range("K19").Value = 275.16666
$ws.Range("L19").Value = 414.57144
$ws.Range("M19").Value = -100.16666
$ws.Range("N19").Value = -764.5714399999999
$ws.Range("H68").Value = 33000
$ws.Range("J68").Value = 33000
$ws.Range("L68").Value = 33000
$ws.Range("N68").Value = -34498
$ws.Range("H71").Value = 33000
$ws.Range("J71").Value = 33000
$ws.Range("L71").Value = 99000
$ws.Range("N71").Value = -106488
$ws.Range("H82").Value = 200
$ws.Range("I82").Value = 200
$ws.Range("K82").Value = 600
$ws.Range("M82").Value = -194
$ws.Range("H85").Value = 200
$ws.Range("I85").Value = 200
$ws.Range("K85").Value = 600
$ws.Range("M85").Value = 804
$ws.Range("H96").Value = 2294.889
$ws.Range("I96").Value = 2941.1667
$ws.Range("J96").Value = 1002.3333
$ws.Range("K96").Value = 8823.500100000001
$ws.Range("L96").Value = 3006.9999
$ws.Range("M96").Value = -7450.500100000001
$ws.Range("N96").Value = -5752.9999
$ws.Range("H97").Value = 500
$ws.Range("J97").Value = 500
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
$ws.Range("H100").Value = 900
$ws.Range("I100").Value = 720
$ws.Range("J100").Value = 2700
$ws.Range("K100").Value = 720
$ws.Range("L100").Value = 2700
$ws.Range("M100").Value = -179
$ws.Range("N100").Value = -3782

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 803.5
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -384
$ws.Range("H74").Value = 1518.8948
$ws.Range("I74").Value = 1063.9333
$ws.Range("J74").Value = 3225
$ws.Range("K74").Value = 1063.9333
$ws.Range("L74").Value = 3225
$ws.Range("M74").Value = -189.9332999999999
$ws.Range("N74").Value = -4973
$ws.Range("H77").Value = 1518.8948
$ws.Range("I77").Value = 1063.9333
$ws.Range("J77").Value = 3225
$ws.Range("K77").Value = 5319.666499999999
$ws.Range("L77").Value = 16125
$ws.Range("M77").Value = -951.6664999999994
$ws.Range("N77").Value = -24861
$ws.Range("H102").Value = 10418057
$ws.Range("I102").Value = 11905923
$ws.Range("K102").Value = 11905923
$ws.Range("M102").Value = -11904301
$ws.Range("H132").Value = 2322.5898
$ws.Range("I132").Value = 1935
$ws.Range("K132").Value = 5805
$ws.Range("M132").Value = -3275

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 125001280
$ws.Range("I99").Value = 166667800
$ws.Range("K99").Value = 166667800
$ws.Range("M99").Value = -166666302

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 176.25
$ws.Range("I7").Value = 74.333336
$ws.Range("J7").Value = 237.4
$ws.Range("K7").Value = 74.333336
$ws.Range("L7").Value = 237.4
$ws.Range("M7").Value = 38.666664
$ws.Range("N7").Value = -463.4
$ws.Range("H31").Value = 1536.0731
$ws.Range("I31").Value = 1447.6666
$ws.Range("K31").Value = 1447.6666
$ws.Range("M31").Value = -1152.6666
$ws.Range("H34").Value = 1536.0731
$ws.Range("I34").Value = 1447.6666
$ws.Range("K34").Value = 1447.6666
$ws.Range("M34").Value = -1245.6666
$ws.Range("H86").Value = 4803537
$ws.Range("I86").Value = 9561095
$ws.Range("J86").Value = 45979.57
$ws.Range("K86").Value = 9561095
$ws.Range("L86").Value = 45979.57
$ws.Range("M86").Value = -9559972
$ws.Range("N86").Value = -48225.57
$ws.Range("H89").Value = 4803537
$ws.Range("I89").Value = 9561095
$ws.Range("J89").Value = 45979.57
$ws.Range("K89").Value = 47805475
$ws.Range("L89").Value = 229897.85
$ws.Range("M89").Value = -47799859
$ws.Range("N89").Value = -241129.85
$ws.Range("H94").Value = 2276.5
$ws.Range("I94").Value = 1703
$ws.Range("K94").Value = 1703
$ws.Range("M94").Value = -1252
$ws.Range("H122").Value = 1006.3333
$ws.Range("I122").Value = 916
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2748
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -298
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 125.55556
$ws.Range("J12").Value = 99.833336
$ws.Range("L12").Value = 299.500008
$ws.Range("N12").Value = -645.500008
$ws.Range("H131").Value = 13334209
$ws.Range("J131").Value = 990.7414
$ws.Range("L131").Value = 2972.2242
$ws.Range("N131").Value = -13052.2242

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5203.4
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10302
$ws.Range("H46").Value = 10250
$ws.Range("I46").Value = 2333.3333
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 2333.3333
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -2177.3333
$ws.Range("N46").Value = -15312
$ws.Range("H48").Value = 9000
$ws.Range("J48").Value = 9000
$ws.Range("L48").Value = 9000
$ws.Range("N48").Value = -9970
$ws.Range("H49").Value = 17250
$ws.Range("J49").Value = 17250
$ws.Range("L49").Value = 17250
$ws.Range("N49").Value = -17618
$ws.Range("H97").Value = 621.17645
$ws.Range("I97").Value = 602.9167
$ws.Range("K97").Value = 602.9167
$ws.Range("M97").Value = -106.9167
$ws.Range("H102").Value = 1177.8889
$ws.Range("I102").Value = 1294.8889
$ws.Range("J102").Value = 1060.8889
$ws.Range("K102").Value = 1294.8889
$ws.Range("L102").Value = 1060.8889
$ws.Range("M102").Value = 327.1111000000001
$ws.Range("N102").Value = -4304.8889
$ws.Range("H126").Value = 1634
$ws.Range("I126").Value = 1252.25
$ws.Range("J126").Value = 2397.5
$ws.Range("K126").Value = 3756.75
$ws.Range("L126").Value = 7192.5
$ws.Range("M126").Value = -1286.75
$ws.Range("N126").Value = -12132.5
$ws.Range("H132").Value = 3519.7727
$ws.Range("I132").Value = 3208.2
$ws.Range("J132").Value = 4187.4287
$ws.Range("K132").Value = 9624.599999999999
$ws.Range("L132").Value = 12562.2861
$ws.Range("M132").Value = -7094.599999999999
$ws.Range("N132").Value = -17622.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2052.7058
$ws.Range("I7").Value = 1947.9231
$ws.Range("J7").Value = 2393.25
$ws.Range("K7").Value = 1947.9231
$ws.Range("L7").Value = 2393.25
$ws.Range("M7").Value = -1835.9231
$ws.Range("N7").Value = -2617.25
$ws.Range("H29").Value = 900
$ws.Range("J29").Value = 900
$ws.Range("L29").Value = 900
$ws.Range("N29").Value = -1490
$ws.Range("H32").Value = 3875
$ws.Range("I32").Value = 2500
$ws.Range("K32").Value = 2500
$ws.Range("M32").Value = -2183
$ws.Range("H34").Value = 2333.3333
$ws.Range("J34").Value = 2500
$ws.Range("L34").Value = 2500
$ws.Range("N34").Value = -2844
$ws.Range("H35").Value = 1950
$ws.Range("J35").Value = 2000
$ws.Range("L35").Value = 2000
$ws.Range("N35").Value = -2672
$ws.Range("H42").Value = 18000
$ws.Range("J42").Value = 18000
$ws.Range("L42").Value = 18000
$ws.Range("N42").Value = -19126
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386
$ws.Range("H49").Value = 18000
$ws.Range("J49").Value = 18000
$ws.Range("L49").Value = 18000
$ws.Range("N49").Value = -18294
$ws.Range("H126").Value = 2052.7058
$ws.Range("I126").Value = 1947.9231
$ws.Range("J126").Value = 2393.25
$ws.Range("K126").Value = 5843.7693
$ws.Range("L126").Value = 7179.75
$ws.Range("M126").Value = -3373.7693
$ws.Range("N126").Value = -12119.75
$ws.Range("H136").Value = 1557.1428
$ws.Range("I136").Value = 1316.6666
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3949.9998
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1399.9998
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H100").Value = 266.8095
$ws.Range("I100").Value = 256.92856
$ws.Range("K100").Value = 513.85712
$ws.Range("M100").Value = 27.14287999999999
$ws.Range("H126").Value = 83338936
$ws.Range("I126").Value = 111114460
$ws.Range("K126").Value = 333343380
$ws.Range("M126").Value = -333340910
